$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 134 - this shifts rows 134..150
# down to 135..151 (and the sheet dimension grows to A1:R151).
$ws.Rows(134).Insert()

# Populate the newly inserted row 134 with the new weekly price record.
$ws.Cells.Item(134, 1).Value  = 7
$ws.Cells.Item(134, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(134, 3).Value  = "Ñuble"
$ws.Cells.Item(134, 4).Value  = 44522
$ws.Cells.Item(134, 5).Value  = 16
$ws.Cells.Item(134, 6).Value  = 100112003
$ws.Cells.Item(134, 7).Value  = "Ajo"
$ws.Cells.Item(134, 8).Value  = "Chino"
$ws.Cells.Item(134, 9).Value  = "Primera"
$ws.Cells.Item(134, 10).Value = 100
$ws.Cells.Item(134, 11).Value = 17000
$ws.Cells.Item(134, 12).Value = 18000
$ws.Cells.Item(134, 13).Value = 17500
$ws.Cells.Item(134, 14).Value = '$/caja 10 kilos'
$ws.Cells.Item(134, 15).Value = "China"
$ws.Cells.Item(134, 16).Value = 1750
$ws.Cells.Item(134, 17).Value = 10
$ws.Cells.Item(134, 18).Value = "Hortaliza"
